# The document contains a single table with 5 "data" rows (rows 1, 5, 9,
# 13, 17 - the others are blank spacer rows) of 5 division problems each.
# Several of the original expressions (e.g. "63÷2=", "10÷6=") repeat more
# than once in the document, and a plain Find/Replace("...", ReplaceAll)
# would rewrite every occurrence identically. Since each occurrence needs
# a different replacement value, address each cell directly instead.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "20÷8="
$t.Cell(1, 2).Range.Text = "88÷9="
$t.Cell(1, 3).Range.Text = "93÷5="
$t.Cell(1, 4).Range.Text = "72÷6="
$t.Cell(1, 5).Range.Text = "35÷6="

$t.Cell(5, 1).Range.Text = "95÷9="
$t.Cell(5, 2).Range.Text = "79÷2="
$t.Cell(5, 3).Range.Text = "13÷8="
$t.Cell(5, 4).Range.Text = "34÷4="
$t.Cell(5, 5).Range.Text = "26÷4="

$t.Cell(9, 1).Range.Text = "56÷5="
$t.Cell(9, 2).Range.Text = "62÷9="
$t.Cell(9, 3).Range.Text = "34÷7="
$t.Cell(9, 4).Range.Text = "89÷6="
$t.Cell(9, 5).Range.Text = "80÷7="

$t.Cell(13, 1).Range.Text = "83÷6="
$t.Cell(13, 2).Range.Text = "80÷2="
$t.Cell(13, 3).Range.Text = "26÷9="
$t.Cell(13, 4).Range.Text = "80÷6="
$t.Cell(13, 5).Range.Text = "39÷9="

$t.Cell(17, 1).Range.Text = "65÷6="
$t.Cell(17, 2).Range.Text = "19÷4="
$t.Cell(17, 3).Range.Text = "12÷5="
$t.Cell(17, 4).Range.Text = "97÷4="
$t.Cell(17, 5).Range.Text = "30÷6="
